# Populate the (previously empty) worksheet with the "failed parses" table
# and style the header row (bold, thin box border, centered / top aligned).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row + 3 data rows, 11 columns (A:K)
$values = New-Object 'object[,]' 4,11

# Row 1 - header
$values[0,0]  = "title"
$values[0,1]  = "text"
$values[0,2]  = "publish_date"
$values[0,3]  = "publish_date_source"
$values[0,4]  = "authors"
$values[0,5]  = "canonical_link"
$values[0,6]  = "feed_link"
$values[0,7]  = "media_link"
$values[0,8]  = "media_title"
$values[0,9]  = "exception_class"
$values[0,10] = "exception_text"

# Row 2
$values[1,0]  = "Stock Market Today: Dow, S&P Live Updates for Oct. 3 - Bloomberg"
$values[1,1]  = ""
$values[1,2]  = "2023-10-02 22:10:28"
$values[1,3]  = "approximated"
$values[1,4]  = ""
$values[1,5]  = ""
$values[1,6]  = "https://news.google.com/rss/articles/CBMicWh0dHBzOi8vd3d3LmJsb29tYmVyZy5jb20vbmV3cy9hcnRpY2xlcy8yMDIzLTEwLTAyL2FzaWEtc3RvY2tzLXNldC10by1kcm9wLWFmdGVyLXNsdW1wLWluLXRyZWFzdXJpZXMtbWFya2V0cy13cmFw0gEA?oc=5"
$values[1,7]  = "https://www.bloomberg.com"
$values[1,8]  = "Bloomberg"
$values[1,9]  = "ValueError"
$values[1,10] = "Essential fields are empty, possibly due to bot protection or bad parse"

# Row 3
$values[2,0]  = "Tesla Sales Drop Allows BYD to Close In - Bloomberg Television"
$values[2,1]  = ""
$values[2,2]  = "2023-10-03 03:24:05"
$values[2,3]  = "approximated"
$values[2,4]  = ""
$values[2,5]  = ""
$values[2,6]  = "https://news.google.com/rss/articles/CCAiC3ZUVm10M2hGdDVzmAEB?oc=5"
$values[2,7]  = "https://www.youtube.com"
$values[2,8]  = "Bloomberg Television"
$values[2,9]  = "ValueError"
$values[2,10] = "Essential fields are empty, possibly due to bot protection or bad parse"

# Row 4
$values[3,0]  = "Abercrombie & Fitch launches investigation into ex-CEO sexual misconduct claims - BBC"
$values[3,1]  = ""
$values[3,2]  = "2023-10-03 10:50:16"
$values[3,3]  = "approximated"
$values[3,4]  = ""
$values[3,5]  = ""
$values[3,6]  = "https://news.google.com/rss/articles/CBMiJ2h0dHBzOi8vd3d3LmJiYy5jb20vbmV3cy93b3JsZC02Njk5MDYyMtIBK2h0dHBzOi8vd3d3LmJiYy5jb20vbmV3cy93b3JsZC02Njk5MDYyMi5hbXA?oc=5"
$values[3,7]  = "https://www.bbc.com"
$values[3,8]  = "BBC"
$values[3,9]  = "ArticleException"
$values[3,10] = "Article ``download()`` failed with HTTPSConnectionPool(host='www.bbc.com', port=443): Read timed out. (read timeout=7) on URL https://www.bbc.com/news/world-66990622"

$ws.Range("A1:K4").Value = $values

# Style the header row: bold font, thin border all around, centered horizontally,
# top-aligned vertically. Build the style on the first header cell, then copy its
# format across the rest of the header row so every header cell shares one style.
$firstHeaderCell = $ws.Range("A1")
$firstHeaderCell.Font.Bold = $true
$firstHeaderCell.Borders.LineStyle = 1
$firstHeaderCell.HorizontalAlignment = -4108
$firstHeaderCell.VerticalAlignment = -4160

$firstHeaderCell.Copy()
$ws.Range("B1:K1").PasteSpecial(-4122)
